# Update the workers_rank_mat ranking table: the underlying matrices scores
# were recomputed, which re-sorts several tied/near-tied rows (changing their
# index/prolificid/name/race) and updates every "matrices" score value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "60bd88b8fc436774352f53b9"
$ws.Range("E2").Value = "Annes"
$ws.Range("G2").Value = 13.33466479435452
$ws.Range("H2").Value = "Asian"

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("E3").Value = "Colleen"
$ws.Range("G3").Value = 13.28508527389088
$ws.Range("H3").Value = "White"

$ws.Range("C4").Value = 19
$ws.Range("D4").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("E4").Value = "Jewel"
$ws.Range("G4").Value = 8.341848052842526

$ws.Range("C5").Value = 22
$ws.Range("D5").Value = "608b14a312c099ac00b721b6"
$ws.Range("E5").Value = "Khushi"
$ws.Range("G5").Value = 8.243276706149189
$ws.Range("H5").Value = "Asian"

$ws.Range("C6").Value = 21
$ws.Range("D6").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("E6").Value = "Bri"
$ws.Range("G6").Value = 8.088597951143703
$ws.Range("H6").Value = "Black or African American"

$ws.Range("C7").Value = 30
$ws.Range("D7").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("E7").Value = "Shadaisia"
$ws.Range("G7").Value = 5.495565970904949
$ws.Range("H7").Value = "Black or African American"

$ws.Range("C8").Value = 33
$ws.Range("D8").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("E8").Value = "Shaniek"
$ws.Range("G8").Value = 5.313043414706029

$ws.Range("C9").Value = 32
$ws.Range("D9").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("E9").Value = "Kellie"
$ws.Range("G9").Value = 5.270427976346927
$ws.Range("H9").Value = "White"

$ws.Range("G10").Value = 4.427089860512693

$ws.Range("G11").Value = 4.1449432509138

$ws.Range("G12").Value = 2.219960194385904

$ws.Range("G13").Value = 1.256578053837044

$ws.Range("G14").Value = 14.28972879330137

$ws.Range("G15").Value = 13.42737676672315

$ws.Range("G16").Value = 8.447915419544929

$ws.Range("G17").Value = 7.374783189707426

$ws.Range("G18").Value = 6.318117057861989

$ws.Range("G19").Value = 6.269575784930271

$ws.Range("C20").Value = 33
$ws.Range("D20").Value = "60b322994d0b901954690036"
$ws.Range("E20").Value = "Brennan"
$ws.Range("G20").Value = 5.377829192037456
$ws.Range("H20").Value = "White"

$ws.Range("C21").Value = 32
$ws.Range("D21").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("E21").Value = "Jamarii"
$ws.Range("G21").Value = 5.220260424381992
$ws.Range("H21").Value = "Black or African American"

$ws.Range("C22").Value = 30
$ws.Range("D22").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("E22").Value = "Matthew"
$ws.Range("G22").Value = 5.0041670821604

$ws.Range("G23").Value = 3.086411040790361

$ws.Range("G24").Value = 1.224224424130261

$ws.Range("G25").Value = 0.3847272879396543
